$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2026-01-03 Saturday", $false, $false, $false, $false, $false, $true, 1, $false, "2026-01-04 Sunday", 2)

# Update the multiplication answers in the table. Cell addressing is used
# (rather than a document-wide Find/Replace) because several cells share
# identical "a×b=c" text, so each cell must be targeted individually.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "56×91=5096"
$t.Cell(1, 2).Range.Text = "93×45=4185"
$t.Cell(1, 3).Range.Text = "85×28=2380"
$t.Cell(1, 4).Range.Text = "63×96=6048"
$t.Cell(1, 5).Range.Text = "71×32=2272"

$t.Cell(5, 1).Range.Text = "45×87=3915"
$t.Cell(5, 2).Range.Text = "29×49=1421"
$t.Cell(5, 3).Range.Text = "55×20=1100"
$t.Cell(5, 4).Range.Text = "65×34=2210"
$t.Cell(5, 5).Range.Text = "15×88=1320"

$t.Cell(10, 1).Range.Text = "52×53=2756"
$t.Cell(10, 2).Range.Text = "26×49=1274"
$t.Cell(10, 3).Range.Text = "41×22=902"
$t.Cell(10, 4).Range.Text = "75×36=2700"
$t.Cell(10, 5).Range.Text = "56×87=4872"

$t.Cell(15, 1).Range.Text = "58×96=5568"
$t.Cell(15, 2).Range.Text = "13×19=247"
$t.Cell(15, 3).Range.Text = "74×90=6660"
$t.Cell(15, 4).Range.Text = "71×26=1846"
$t.Cell(15, 5).Range.Text = "37×44=1628"

$t.Cell(20, 1).Range.Text = "93×91=8463"
$t.Cell(20, 2).Range.Text = "39×45=1755"
$t.Cell(20, 3).Range.Text = "98×34=3332"
$t.Cell(20, 4).Range.Text = "70×12=840"
$t.Cell(20, 5).Range.Text = "27×48=1296"

$d.Save()
